# Insert a new weekly record at row 29, pushing all existing rows 29-111
# down to 30-112 (handled automatically by the row insert), then populate
# the newly blank row 29 with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(29).Insert()

$ws.Range("A29").Value = 10
$ws.Range("B29").Value = "Vega Modelo de Temuco"
$ws.Range("C29").Value = "La Araucanía"
$ws.Range("D29").Value = 44925
$ws.Range("E29").Value = 9
$ws.Range("F29").Value = 100112022
$ws.Range("G29").Value = "Arveja Verde"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 120
$ws.Range("K29").Value = 25000
$ws.Range("L29").Value = 25000
$ws.Range("M29").Value = 25000
$ws.Range("N29").Value = "$/saco 25 kilos"
$ws.Range("O29").Value = "Región de La Araucanía"
$ws.Range("P29").Value = 1000
$ws.Range("Q29").Value = 25
$ws.Range("R29").Value = "Hortaliza"
